$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of exported data: a single-cell entry in A2 whose text begins
# with a line break (matches the source export's "\n1224100619Dariel"
# shared string), e.g. an ID/name pair exported together.
$ws.Cells.Item(2, 1).Value = "`n1224100619Dariel"

# Let Excel recompute the row's height from scratch instead of leaving a
# stale autofit/custom height behind from the multi-line text assignment.
$ws.Rows.Item(2).AutoFit()
